$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 2900.5557
$ws.Range("I69").Value = 2000
$ws.Range("J69").Value = 3013.125
$ws.Range("K69").Value = 6000
$ws.Range("L69").Value = 9039.375
$ws.Range("M69").Value = -5126
$ws.Range("N69").Value = -10787.375
$ws.Range("H72").Value = 2900.5557
$ws.Range("I72").Value = 2000
$ws.Range("J72").Value = 3013.125
$ws.Range("K72").Value = 18000
$ws.Range("L72").Value = 27118.125
$ws.Range("M72").Value = -13632
$ws.Range("N72").Value = -35854.125
$ws.Range("H98").Value = 803.94116
$ws.Range("I98").Value = 785.13336
$ws.Range("J98").Value = 945
$ws.Range("K98").Value = 785.13336
$ws.Range("L98").Value = 945
$ws.Range("M98").Value = 712.86664
$ws.Range("N98").Value = -3941
$ws.Range("H122").Value = 803.94116
$ws.Range("I122").Value = 785.13336
$ws.Range("J122").Value = 945
$ws.Range("K122").Value = 2355.40008
$ws.Range("L122").Value = 2835
$ws.Range("M122").Value = 94.59991999999966
$ws.Range("N122").Value = -7735
$ws.Range("H129").Value = 921.125
$ws.Range("I129").Value = 558.5714
$ws.Range("K129").Value = 1675.7142
$ws.Range("M129").Value = 3324.2858
$ws.Range("H137").Value = 3848201.2
$ws.Range("I137").Value = 4349586.5
$ws.Range("K137").Value = 13048759.5
$ws.Range("M137").Value = -13046209.5
$ws.Range("H138").Value = 3147292
$ws.Range("I138").Value = 1013.8947
$ws.Range("J138").Value = 4905506
$ws.Range("K138").Value = 3041.6841
$ws.Range("L138").Value = 14716518
$ws.Range("M138").Value = 2098.3159
$ws.Range("N138").Value = -14726798

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1613.16
$ws.Range("I32").Value = 1336.337
$ws.Range("K32").Value = 1336.337
$ws.Range("M32").Value = -1049.337
$ws.Range("H74").Value = 6302081
$ws.Range("I74").Value = 9655509
$ws.Range("J74").Value = 74285.71000000001
$ws.Range("K74").Value = 9655509
$ws.Range("L74").Value = 74285.71000000001
$ws.Range("M74").Value = -9654635
$ws.Range("N74").Value = -76033.71000000001
$ws.Range("H77").Value = 6302081
$ws.Range("I77").Value = 9655509
$ws.Range("J77").Value = 74285.71000000001
$ws.Range("K77").Value = 48277545
$ws.Range("L77").Value = 371428.55
$ws.Range("M77").Value = -48273177
$ws.Range("N77").Value = -380164.55
$ws.Range("H97").Value = 3473033.2
$ws.Range("I97").Value = 4167363.2
$ws.Range("J97").Value = 1383.3334
$ws.Range("K97").Value = 4167363.2
$ws.Range("L97").Value = 1383.3334
$ws.Range("M97").Value = -4166867.2
$ws.Range("N97").Value = -2375.3334
$ws.Range("H110").Value = 910682.0600000001
$ws.Range("I110").Value = 1250873.8
$ws.Range("J110").Value = 3504.3333
$ws.Range("K110").Value = 1250873.8
$ws.Range("L110").Value = 3504.3333
$ws.Range("M110").Value = -1248828.8
$ws.Range("N110").Value = -7594.3333
$ws.Range("H132").Value = 35613.31
$ws.Range("I132").Value = 25771.8
$ws.Range("J132").Value = 57483.332
$ws.Range("K132").Value = 77315.39999999999
$ws.Range("L132").Value = 172449.996
$ws.Range("M132").Value = -74785.39999999999
$ws.Range("N132").Value = -177509.996

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 20000
$ws.Range("I26").Value = 20000
$ws.Range("K26").Value = 20000
$ws.Range("M26").Value = -19708
$ws.Range("H86").Value = 18808.25
$ws.Range("I86").Value = 15773.588
$ws.Range("J86").Value = 36004.668
$ws.Range("K86").Value = 15773.588
$ws.Range("L86").Value = 36004.668
$ws.Range("M86").Value = -14650.588
$ws.Range("N86").Value = -38250.668
$ws.Range("H89").Value = 18808.25
$ws.Range("I89").Value = 15773.588
$ws.Range("J89").Value = 36004.668
$ws.Range("K89").Value = 78867.94
$ws.Range("L89").Value = 180023.34
$ws.Range("M89").Value = -73251.94
$ws.Range("N89").Value = -191255.34
$ws.Range("H94").Value = 467.7647
$ws.Range("I94").Value = 452.66666
$ws.Range("J94").Value = 504
$ws.Range("K94").Value = 452.66666
$ws.Range("L94").Value = 504
$ws.Range("M94").Value = -1.666659999999979
$ws.Range("N94").Value = -1406
$ws.Range("H107").Value = 1815.1666
$ws.Range("I107").Value = 1871
$ws.Range("K107").Value = 1871
$ws.Range("M107").Value = 49
$ws.Range("H134").Value = 2422.2263
$ws.Range("I134").Value = 1897.2325
$ws.Range("K134").Value = 5691.6975
$ws.Range("M134").Value = -3156.6975

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H18").Value = 40684
$ws.Range("J18").Value = 40684
$ws.Range("L18").Value = 40684
$ws.Range("H114").Value = 49970
$ws.Range("J114").Value = 49970
$ws.Range("L114").Value = 49970
$ws.Range("H132").Value = 22544.104
$ws.Range("I132").Value = 1376.4878
$ws.Range("J132").Value = 146525.86
$ws.Range("K132").Value = 4129.463400000001
$ws.Range("L132").Value = 439577.58
$ws.Range("M132").Value = -1599.463400000001
$ws.Range("N132").Value = -444637.58

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 848.0476
$ws.Range("I5").Value = 433.16666
$ws.Range("K5").Value = 1299.49998
$ws.Range("M5").Value = -1187.49998
$ws.Range("H113").Value = 679.36
$ws.Range("I113").Value = 563.5
$ws.Range("J113").Value = 701.4286
$ws.Range("K113").Value = 1690.5
$ws.Range("L113").Value = 2104.2858
$ws.Range("M113").Value = 479.5
$ws.Range("N113").Value = -6444.2858
$ws.Range("H131").Value = 1014.9359
$ws.Range("I131").Value = 696.6667
$ws.Range("J131").Value = 1041.4584
$ws.Range("K131").Value = 2090.0001
$ws.Range("L131").Value = 3124.3752
$ws.Range("M131").Value = 2949.9999
$ws.Range("N131").Value = -13204.3752
$ws.Range("H135").Value = 848.0476
$ws.Range("I135").Value = 433.16666
$ws.Range("K135").Value = 3898.49994
$ws.Range("M135").Value = -1363.49994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3797.2104
$ws.Range("I80").Value = 3240
$ws.Range("K80").Value = 3240
$ws.Range("M80").Value = -2242
$ws.Range("H83").Value = 3797.2104
$ws.Range("I83").Value = 3240
$ws.Range("K83").Value = 16200
$ws.Range("M83").Value = -11208
$ws.Range("H97").Value = 1451.4615
$ws.Range("I97").Value = 1811.125
$ws.Range("J97").Value = 876
$ws.Range("K97").Value = 1811.125
$ws.Range("L97").Value = 876
$ws.Range("M97").Value = -1315.125
$ws.Range("H122").Value = 2439.6667
$ws.Range("I122").Value = 1979
$ws.Range("J122").Value = 4052
$ws.Range("K122").Value = 5937
$ws.Range("L122").Value = 12156
$ws.Range("M122").Value = -3487
$ws.Range("N122").Value = -17056
$ws.Range("H126").Value = 1591.6666
$ws.Range("I126").Value = 1410
$ws.Range("J126").Value = 2500
$ws.Range("K126").Value = 4230
$ws.Range("L126").Value = 7500
$ws.Range("M126").Value = -1760
$ws.Range("N126").Value = -12440
$ws.Range("H132").Value = 48014.465
$ws.Range("I132").Value = 31548.182
$ws.Range("J132").Value = 102353.2
$ws.Range("K132").Value = 94644.546
$ws.Range("L132").Value = 307059.6
$ws.Range("M132").Value = -92114.546
$ws.Range("N132").Value = -312119.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1403.2
$ws.Range("I100").Value = 1171.1111
$ws.Range("J100").Value = 2000
$ws.Range("K100").Value = 1171.1111
$ws.Range("L100").Value = 2000
$ws.Range("M100").Value = -630.1111000000001
$ws.Range("N100").Value = -3082
$ws.Range("H132").Value = 36432.16
$ws.Range("I132").Value = 15263.605
$ws.Range("J132").Value = 170499.67
$ws.Range("K132").Value = 45790.815
$ws.Range("L132").Value = 511499.01
$ws.Range("M132").Value = -43260.815
$ws.Range("N132").Value = -516559.01

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 385.06668
$ws.Range("I107").Value = 290
$ws.Range("J107").Value = 468.25
$ws.Range("K107").Value = 870
$ws.Range("L107").Value = 1404.75
$ws.Range("M107").Value = 1050
$ws.Range("N107").Value = -5244.75
$ws.Range("H132").Value = 51649
$ws.Range("I132").Value = 51121.15
$ws.Range("J132").Value = 52176.85
$ws.Range("K132").Value = 153363.45
$ws.Range("L132").Value = 156530.55
$ws.Range("M132").Value = -150833.45
$ws.Range("N132").Value = -161590.55
$ws.Range("H133").Value = 40150
$ws.Range("J133").Value = 40150
$ws.Range("L133").Value = 40150
$ws.Range("N133").Value = -50270
$ws.Range("H136").Value = 38336.02
$ws.Range("I136").Value = 20982.9
$ws.Range("K136").Value = 62948.7
$ws.Range("M136").Value = -70850
